$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A2" = "5.175744445296004e-05"
    "B2" = "0.0002706766244955361"
    "C2" = "2.867047987820115e-05"
    "D2" = "0.0001835719886003062"
    "E2" = "6.814509106334299e-05"
    "F2" = "9.353094355901703e-05"
    "G2" = "3.372250648681074e-05"
    "H2" = "5.803597377962433e-05"
    "I2" = "4.642974090529606e-06"
    "J2" = "1.882780634332448e-05"
    "K2" = "0.0002689994289539754"
    "L2" = "2.096350362990052e-05"
    "M2" = "0.0002048777241725475"
    "N2" = "3.265701525378972e-05"
    "O2" = "0.0001726726914057508"
    "P2" = "8.53581732371822e-05"
    "Q2" = "0.0001289033680222929"
    "R2" = "0.0001509142894065008"
    "S2" = "1.90005375770852e-05"
    "T2" = "5.22313712281175e-05"
    "U2" = "1.077256729331566e-05"
    "V2" = "8.388151400140487e-06"
    "W2" = "4.74683110951446e-05"
    "X2" = "1.704878741293214e-06"
    "Y2" = "7.007025851635262e-05"
    "Z2" = "2.18781060539186e-05"
    "AA2" = "3.614824890973978e-05"
    "AB2" = "0.0001250212080776691"
    "AC2" = "1.613634231034666e-05"
    "AD2" = "6.549978934344836e-06"
    "AE2" = "6.453505193348974e-05"
    "AF2" = "4.138883377891034e-05"
    "AG2" = "1.942700691870414e-05"
    "AH2" = "2.121151010214817e-05"
    "AI2" = "2.991000837937463e-05"
    "AJ2" = "5.934435466770083e-05"
    "AK2" = "7.393797568511218e-05"
    "AL2" = "0.000124809259432368"
    "AM2" = "1.004185742203845e-05"
    "AN2" = "6.31367220194079e-05"
    "AO2" = "4.590603930409998e-05"
    "AP2" = "2.27569626076729e-06"
    "AQ2" = "5.760020940215327e-06"
    "AR2" = "6.501402822323143e-05"
    "AS2" = "5.738116669817828e-06"
    "AT2" = "0.0001083151801140048"
    "AU2" = "0.0001218029428855516"
    "AV2" = "1.304979286942398e-05"
    "AW2" = "4.323662506067194e-05"
    "AX2" = "1.21126704470953e-05"
    "AY2" = "4.521081427810714e-05"
    "AZ2" = "7.040426135063171e-05"
    "BA2" = "2.102085272781551e-05"
    "BB2" = "1.034011620504316e-05"
    "BC2" = "3.209827809769195e-06"
    "BD2" = "0.0001403368805767968"
    "BE2" = "1.266394247068092e-05"
    "BF2" = "9.164206858258694e-05"
    "BG2" = "5.62178720429074e-05"
    "BH2" = "1.34830224851612e-05"
    "BI2" = "1.176955083792564e-05"
    "BJ2" = "2.180605588364415e-05"
    "BK2" = "3.624587407102808e-05"
    "BL2" = "9.270802547689527e-05"
    "BM2" = "0.0001286218903260306"
    "BN2" = "1.499241716373945e-05"
    "BO2" = "7.488504343200475e-05"
    "BP2" = "7.008493412286043e-05"
    "BQ2" = "1.7034250049619e-05"
    "BR2" = "1.018118928186595e-05"
    "BS2" = "3.108203964075074e-05"
    "BT2" = "3.764618668355979e-05"
    "BU2" = "8.880347013473511e-05"
    "BV2" = "8.856803469825536e-05"
    "BW2" = "1.857190727605484e-05"
    "BX2" = "3.038812792510726e-06"
    "BY2" = "0.000106357409094926"
    "BZ2" = "3.230220681871288e-05"
    "CA2" = "5.190101364860311e-05"
    "CB2" = "3.778905374929309e-05"
    "CC2" = "7.108539284672588e-05"
    "CD2" = "5.404969124356285e-05"
    "CE2" = "1.757677091518417e-05"
    "CF2" = "9.599605618859641e-06"
    "CG2" = "5.347168917069212e-05"
    "CH2" = "5.532753857551143e-05"
    "CI2" = "6.687891072942875e-06"
    "CJ2" = "1.627500023460016e-05"
    "CK2" = "6.253741503314814e-06"
    "CL2" = "5.935498393228045e-06"
    "CM2" = "3.457839557086118e-05"
    "CN2" = "2.23466176976217e-05"
    "CO2" = "1.170488940260839e-05"
    "CP2" = "3.669754732982256e-05"
    "CQ2" = "1.265115315618459e-05"
    "CR2" = "8.99236329132691e-05"
    "CS2" = "0.0001064281605067663"
    "CT2" = "0.0001549352746224031"
    "CU2" = "9.875572868622839e-05"
    "CV2" = "5.107189645059407e-05"
    "CW2" = "1.716981387289707e-05"
    "CX2" = "1.802030783437658e-05"
    "CY2" = "4.112636452191509e-05"
    "CZ2" = "3.15225770464167e-05"
    "DA2" = "1.940890615514945e-05"
    "DB2" = "3.267823558417149e-06"
    "DC2" = "7.363026725215605e-06"
    "DD2" = "3.69497547580977e-06"
    "DE2" = "3.004015707119834e-05"
    "DF2" = "0.0002082170103676617"
    "DG2" = "8.101592538878322e-05"
    "DH2" = "0.0001255649694940075"
    "DI2" = "6.136229785624892e-05"
    "DJ2" = "0.0001056004839483649"
    "DK2" = "0.0001508198329247534"
    "DL2" = "1.999964842980262e-05"
    "DM2" = "7.056284812279046e-05"
    "DN2" = "9.956368012353778e-05"
    "DO2" = "6.681309605482966e-05"
    "DP2" = "0.0001330270170001313"
    "DQ2" = "0.000121326265798416"
    "DR2" = "5.025275459047407e-05"
    "DS2" = "3.294432099210098e-05"
    "DT2" = "1.181156767415814e-05"
    "DU2" = "5.418379805632867e-05"
    "DV2" = "2.132661211362574e-05"
    "DW2" = "9.479768050368875e-05"
    "DX2" = "5.807372508570552e-05"
    "DY2" = "3.212837327737361e-05"
    "DZ2" = "4.564048504107632e-05"
    "EA2" = "1.262151272385381e-05"
    "EB2" = "0.000130759333842434"
    "EC2" = "6.293338174145902e-06"
    "ED2" = "4.720764263765886e-05"
    "EE2" = "5.252527262200601e-05"
    "EF2" = "5.768742994405329e-05"
    "EG2" = "2.519247391319368e-05"
    "EH2" = "4.47930688096676e-08"
    "EI2" = "8.955871453508735e-05"
    "EJ2" = "1.581421383889392e-05"
    "EK2" = "2.367903471167665e-05"
    "EL2" = "2.59720764006488e-07"
    "EM2" = "3.23223284794949e-05"
    "EN2" = "3.327148533571744e-06"
    "EO2" = "6.509016384370625e-05"
    "EP2" = "2.528537879697978e-05"
    "EQ2" = "9.402963041793555e-05"
    "ER2" = "5.053755376138724e-05"
    "ES2" = "8.335719030583277e-06"
    "ET2" = "7.895525050116703e-05"
    "EU2" = "1.505700674897525e-05"
    "EV2" = "6.620703970838804e-06"
    "EW2" = "6.780916010029614e-06"
    "EX2" = "6.302439578576013e-05"
    "EY2" = "1.900684037536848e-05"
    "EZ2" = "2.051655792456586e-05"
    "FA2" = "5.231084287515841e-05"
    "FB2" = "4.479569179238752e-05"
    "FC2" = "2.933199766630423e-06"
    "FD2" = "1.698635605862364e-05"
    "FE2" = "5.976000920782099e-06"
    "FF2" = "1.334234002570156e-05"
    "FG2" = "1.294805588258896e-05"
    "FH2" = "9.903916361508891e-05"
    "FI2" = "3.001097866217606e-05"
    "FJ2" = "6.301519169937819e-05"
    "FK2" = "9.975024113373365e-06"
    "FL2" = "2.250008401460946e-05"
    "FM2" = "2.655043499544263e-05"
    "FN2" = "2.586860864539631e-05"
    "FO2" = "4.394241113914177e-05"
    "FP2" = "0.0001498667697887868"
    "FQ2" = "0.000100951649073977"
    "FR2" = "4.86735807498917e-05"
    "FS2" = "8.863559196470305e-05"
    "FT2" = "2.088681321765762e-05"
    "FU2" = "4.396339136292227e-06"
    "FV2" = "0.0001755518896970898"
    "FW2" = "9.574603609507903e-05"
    "FX2" = "4.7238492697943e-05"
    "FY2" = "0.0001513330062152818"
    "FZ2" = "3.979960638389457e-06"
    "GA2" = "0.0001558093354105949"
    "GB2" = "1.097512722481042e-05"
    "GC2" = "0.0001049283164320514"
    "GD2" = "7.364054181380197e-05"
    "GE2" = "6.014137034071609e-05"
    "GF2" = "2.407970896456391e-06"
    "GG2" = "2.080161357298493e-05"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = [double]$updates[$addr]
}
